$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: replace "z" with "b" in B2:E2
$ws.Range("B2").Value = "b"
$ws.Range("C2").Value = "b"
$ws.Range("D2").Value = "b"
$ws.Range("E2").Value = "b"

# Add new row 4 with new data
$ws.Range("A4").Value = "la_bagguette"
$ws.Range("B4").Value = "Piastri"
$ws.Range("C4").Value = "Gasly"
$ws.Range("D4").Value = "Ocon"
$ws.Range("E4").Value = "Verstappen"
